$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (column D) and Volume(1h) (column E) values.
# Column D values are forced to Text so numeric-looking strings (e.g. "582.36")
# are preserved exactly as text instead of being coerced into floating point
# numbers; ClearFormats() afterwards removes the temporary "@" number format
# so the cell keeps its original (default) style, matching the source file.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.122.68'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.468.29'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.93'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.917.16'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.36'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.047.59'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.415.00'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.90'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.44'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '348.33'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.01'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.18'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.595.15'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0897'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '498.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.73'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.82'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.14'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.80'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '142.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₆0255'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.508'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('E51').Value = '  -0.36%  '
